$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "Seqance" -> "Sequence" in D1
$ws.Range("D1").Value = "Sequence"

# Update the selection to D1
$ws.Range("D1").Select()
